# Insert a new weekly price record above the current row 478 on the
# "Feria Lagunitas de Puerto Montt - Brócoli" sheet. This shifts the
# existing rows 478-498 down to 479-499 (dimension grows to A1:R499)
# and fills the newly opened row 478 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 478..498 down by one to make room for the new record.
$ws.Rows.Item(478).Insert()

# Populate the new row 478 with the new weekly data point.
$ws.Cells.Item(478, 1).Value = 4
$ws.Cells.Item(478, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(478, 3).Value = "Los Lagos"
$ws.Cells.Item(478, 4).Value = 45008
$ws.Cells.Item(478, 5).Value = 10
$ws.Cells.Item(478, 6).Value = 100112023
$ws.Cells.Item(478, 7).Value = "Brócoli"
$ws.Cells.Item(478, 8).Value = "Sin especificar"
$ws.Cells.Item(478, 9).Value = "Primera"
$ws.Cells.Item(478, 10).Value = 500
$ws.Cells.Item(478, 11).Value = 1500
$ws.Cells.Item(478, 12).Value = 1600
$ws.Cells.Item(478, 13).Value = 1550
$ws.Cells.Item(478, 14).Value = "$/unidad"
$ws.Cells.Item(478, 15).Value = "Región Metropolitana"
$ws.Cells.Item(478, 16).Value = 1550
$ws.Cells.Item(478, 17).Value = 1
$ws.Cells.Item(478, 18).Value = "Hortaliza"
